$wb = $excel.ActiveWorkbook

# --- StatOutput_Message: A18 gets the new Cypher query text (Bernese Mountain Dog) ---
$statMessage = $wb.Worksheets.Item("StatOutput_Message")
$a18 = $statMessage.Range("A18")
$a18.Value2 = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bernese Mountain Dog']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# --- StatOutput: C2 gets the new file/case count ("6"), D2 gets the study count ("2") ---
$statOutput = $wb.Worksheets.Item("StatOutput")

$c2 = $statOutput.Range("C2")
$c2.NumberFormat = "@"
$c2.Value2 = "6"

$d2 = $statOutput.Range("D2")
$d2.NumberFormat = "@"
$d2.Value2 = "2"
